# Update Rizka - Data Binding Lihat Status Pengajuan Data Ketinggalan
#
# Sheet1 data table changes:
#   - Row 3 car name "Calya" -> "Kijang Innova"
#   - Row 4 ("Daihatsu New Ayla" / "searchByText") is removed entirely,
#     shifting the last row ("Brio" / "failed") up from row 5 to row 4.
#   - Active selection moves from D2 to A4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the car in row 3 from "Calya" to "Kijang Innova"
$ws.Range("A3").Value = "Kijang Innova"

# Delete row 4 ("Daihatsu New Ayla" / "searchByText"); subsequent rows shift up
$ws.Rows.Item(4).Delete()

# Match the final active selection in the saved file
[void]$ws.Range("A4").Select()
